# Add data for 2022-08-04 (i.e. updates "through" date from 07-26 to 07-27)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab to reflect the new "through" date
$ws.Name = "Through 2022-07-27"

# Update the row label for July
$ws.Range("A8").Value = "July (through 07-27)"

# Updated July counts (row 8), for years 2015-2022 in columns B-I
$julyValues = @(36, 47, 64, 63, 43, 124, 133, 151)
for ($i = 0; $i -lt $julyValues.Length; $i++) {
    $col = 2 + $i
    $ws.Cells.Item(8, $col).Value = $julyValues[$i]
}

# Updated Total counts (row 9), for years 2015-2022 in columns B-I
$totalValues = @(161, 295, 454, 416, 294, 596, 893, 957)
for ($i = 0; $i -lt $totalValues.Length; $i++) {
    $col = 2 + $i
    $ws.Cells.Item(9, $col).Value = $totalValues[$i]
}
